# Hortaliza, Vega Modelo de Temuco - Albahaca
# A new weekly price observation is inserted as row 114, pushing the
# existing rows 114..238 down to 115..239 (dimension grows from
# A1:R238 to A1:R239).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 114, shifting everything
# below it down by one row.
$ws.Rows("114:114").Insert()

# Populate the newly inserted row 114 with the new observation.
$ws.Range("A114").Value = 10
$ws.Range("B114").Value = 'Vega Modelo de Temuco'
$ws.Range("C114").Value = 'La Araucanía'
$ws.Range("D114").Value = 44740
$ws.Range("E114").Value = 9
$ws.Range("F114").Value = 100112052
$ws.Range("G114").Value = 'Albahaca'
$ws.Range("H114").Value = 'Sin especificar'
$ws.Range("I114").Value = 'Primera'
$ws.Range("J114").Value = 20
$ws.Range("K114").Value = 5500
$ws.Range("L114").Value = 5500
$ws.Range("M114").Value = 5500
$ws.Range("N114").Value = '$/paquete'
$ws.Range("O114").Value = 'Región de Arica y Parinacota'
$ws.Range("P114").Value = 5500
$ws.Range("Q114").Value = 1
$ws.Range("R114").Value = 'Hortaliza'

# Match the date number format used by the rest of column D.
$ws.Range("D114").NumberFormat = $ws.Range("D115").NumberFormat
